$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.505.59"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "1.845.20"
$ws.Range("E3").Value = "  +4.07%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.35"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.29%  "
$ws.Range("E9").Value = "  +6.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0719"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0935"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "2.108.95"
$ws.Range("E12").Value = "  +4.01%  "
$ws.Range("D13").Value = "1.845.50"
$ws.Range("E13").Value = "  +4.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.653"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.15%  "
$ws.Range("D16").Value = "34.529.21"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("D20").Value = "0.0₃0806"
$ws.Range("E20").Value = "  +9.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("E23").Value = "  +3.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.88%  "
$ws.Range("E26").Value = "  +3.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.23%  "
$ws.Range("E28").Value = "  +2.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  +5.95%  "
$ws.Range("B31").Value = "Swop.fi"
$ws.Range("C31").Value = "https://coinranking.com/coin/yrCr2HW2c+swopfi-swop"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "535.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +926.63%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.85%  "
$ws.Range("E33").Value = "  +1.96%  "
$ws.Range("E34").Value = "  +2.81%  "
$ws.Range("E35").Value = "  +7.17%  "
$ws.Range("D36").Value = "1.465.00"
$ws.Range("E36").Value = "  +1.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.657"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.74%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0195"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.78%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.978"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "83.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.18%  "
$ws.Range("D46").Value = "2.006.23"
$ws.Range("E46").Value = "  +4.15%  "
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("E48").Value = "  -2.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "106.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.73%  "
$ws.Range("E51").Value = "  +0.13%  "
